$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "brand_id" column (K) is being dropped; the "expertise_id" column
# (currently L) takes its place in column K, one column to the left.
$expertiseHeader = $ws.Range("L1").Value2
$expertiseValue  = $ws.Range("L2").Value2

# Move the expertise_id header/value left into column K (overwriting
# brand_id), then clear the now-duplicate column L.
$ws.Range("K1").Value2 = $expertiseHeader
$ws.Range("K2").Value2 = $expertiseValue
$ws.Range("L1").ClearContents()
$ws.Range("L2").ClearContents()

# K1 had the bordered/centered header look inherited from the old
# "brand_id" cell; the relocated "expertise_id" header instead uses the
# plain bold style (no border/center alignment) since it's now the last
# column.
$ws.Range("K1").Style = "Normal"
$ws.Range("K1").Font.Bold = $true

# Reflect where the user left the selection after the edit.
$ws.Range("K2").Select()
